# Apply updated cryptocurrency price/volume data (GitHub Actions symbol-list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) are stored as plain text in this sheet (e.g. "256.88",
# "0.79%"), so a leading "'" is used to force text entry instead of Excel auto-converting
# the look-alike numeric/percent strings to numbers; the style is then reset to Normal so no
# stray quote-prefix formatting is left behind on the cell.

$ws.Range("D2").Value = "'256.99"
$ws.Range("E2").Value = "'0.93%"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'27.15"
$ws.Range("E3").Value = "'-3.43%"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'4.794"
$ws.Range("E4").Value = "'-9.37%"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.05946"
$ws.Range("E5").Value = "'1.62%"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'6.659"
$ws.Range("E6").Value = "'-0.56%"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.8703"
$ws.Range("E7").Value = "'-0.05%"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.9510"
$ws.Range("E8").Value = "'2.91%"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.1404"
$ws.Range("E9").Value = "'-0.65%"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.03815"
$ws.Range("E10").Value = "'10.30%"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07160"
$ws.Range("E11").Value = "'0.81%"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.03203"
$ws.Range("E12").Value = "'0.78%"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.09260"
$ws.Range("E13").Value = "'0.35%"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.001546"
$ws.Range("E14").Value = "'0.29%"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.0006073"
$ws.Range("E15").Value = "'-94.27%"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.006055"
$ws.Range("E16").Value = "'3.94%"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.481"
$ws.Range("E17").Value = "'-0.47%"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'3.194"
$ws.Range("E18").Value = "'-1.21%"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'2.239"
$ws.Range("E19").Value = "'0.76%"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'-1.41%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'0.1307"
$ws.Range("E21").Value = "'-0.68%"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'3.824"
$ws.Range("E22").Value = "'8.53%"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04202"
$ws.Range("E23").Value = "'0.76%"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'2.58%"
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'-0.72%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.004493"
$ws.Range("E26").Value = "'-10.13%"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.0001200"
$ws.Range("E27").Value = "'-0.05%"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'0.0001494"
$ws.Range("E28").Value = "'86.66%"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Style = "Normal"

$ws.Range("D40").Value = "'0.03841"
$ws.Range("E40").Value = "'0.65%"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006250"
$ws.Range("E41").Value = "'19.44%"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1100"
$ws.Range("E42").Value = "'-0.09%"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.002252"
$ws.Range("E43").Value = "'-3.74%"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'4.51%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005503"
$ws.Range("E45").Value = "'5.35%"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'-0.04%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.08857"
$ws.Range("E47").Value = "'-4.80%"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.002383"
$ws.Range("E48").Value = "'10.47%"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'-0.04%"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'-0.04%"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Style = "Normal"
